$d = $word.ActiveDocument

$replacements = @(
    @("80×80=", "49×64="),
    @("72×50=", "23×38="),
    @("32×80=", "40×97="),
    @("16×41=", "58×86="),
    @("63×39=", "75×36="),
    @("25×24=", "64×47="),
    @("50×75=", "18×14="),
    @("46×79=", "89×22="),
    @("79×56=", "54×43="),
    @("97×77=", "89×45="),
    @("19×65=", "31×24="),
    @("28×20=", "33×87="),
    @("63×81=", "71×81="),
    @("20×72=", "36×99="),
    @("84×51=", "49×70="),
    @("41×95=", "32×20="),
    @("38×61=", "62×31="),
    @("38×48=", "63×62="),
    @("74×82=", "12×90="),
    @("90×84=", "25×73="),
    @("78×36=", "31×82="),
    @("11×89=", "98×21="),
    @("96×43=", "66×82="),
    @("77×69=", "55×56="),
    @("28×81=", "31×55=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
